$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 10: fill in the previously-empty result cells (C,D,E,G,H,I)
# Styles are copied from existing rows that already carry the target
# style indices (s="10", s="12", s="11"), then the values/text are
# overwritten.
# -----------------------------------------------------------------
$ws.Range("C9").Copy($ws.Range("C10"))
$ws.Range("C10").Value = 303

$ws.Range("D9").Copy($ws.Range("D10"))
$ws.Range("D10").Value = "watermelon_mt"

$ws.Range("E9").Copy($ws.Range("E10"))
$ws.Range("E10").Value = 379236

$ws.Range("G9").Copy($ws.Range("G10"))
$ws.Range("G10").Value = 254

$ws.Range("H9").Copy($ws.Range("H10"))
$ws.Range("H10").Value = [double]"4E-70"

$ws.Range("I9").Copy($ws.Range("I10"))
$ws.Range("I10").Value = "3/-3"

# Row 11 grew taller (content elsewhere on the sheet forced a re-fit)
$ws.Rows.Item(11).RowHeight = 21

# -----------------------------------------------------------------
# Row 12: fill in the previously-empty result cells (C,D,E,G,H,I)
# -----------------------------------------------------------------
$ws.Range("C9").Copy($ws.Range("C12"))
$ws.Range("C12").Value = 379236

$ws.Range("D9").Copy($ws.Range("D12"))
$ws.Range("D12").Value = "watermelon_mt"

$ws.Range("E9").Copy($ws.Range("E12"))
$ws.Range("E12").Value = 379236

$ws.Range("H5").Copy($ws.Range("G12"))
$ws.Range("G12").Value = 700300

$ws.Range("C9").Copy($ws.Range("H12"))
$ws.Range("H12").Value = 0

$ws.Range("I9").Copy($ws.Range("I12"))
$ws.Range("I12").Value = "Plus/Plus"

# -----------------------------------------------------------------
# Row 13 (new)
# -----------------------------------------------------------------
$ws.Range("C9").Copy($ws.Range("C13"))
$ws.Range("C13").Value = 379236

$ws.Range("D9").Copy($ws.Range("D13"))
$ws.Range("D13").Value = "watermelon_mt"

$ws.Range("E9").Copy($ws.Range("E13"))
$ws.Range("E13").Value = 379236

$ws.Range("H5").Copy($ws.Range("G13"))
$ws.Range("G13").Value = 683900

$ws.Range("C9").Copy($ws.Range("H13"))
$ws.Range("H13").Value = 0

$ws.Range("I9").Copy($ws.Range("I13"))
$ws.Range("I13").Value = "Plus/Plus"

# -----------------------------------------------------------------
# Row 14 (new)
# -----------------------------------------------------------------
$ws.Range("C9").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 379236

$ws.Range("D9").Copy($ws.Range("D14"))
$ws.Range("D14").Value = "watermelon_mt"

$ws.Range("E9").Copy($ws.Range("E14"))
$ws.Range("E14").Value = 379236

$ws.Range("H5").Copy($ws.Range("G14"))
$ws.Range("G14").Value = 399400

$ws.Range("C9").Copy($ws.Range("H14"))
$ws.Range("H14").Value = 0

$ws.Range("I9").Copy($ws.Range("I14"))
$ws.Range("I14").Value = "Plus/Plus"

# Move the active selection to C14, matching the saved cursor position
$null = $ws.Range("C14").Select()
